# Updates the 'Pais' (countries COVID stats) sheet to a newer data pull.
# - Refreshes the 'datos actualizados' timestamp in A1.
# - Refreshes per-country totals (Casos totales/Nuevos casos/Casos activos/
#   Recuperados/Muertes hoy/Muertes) for the rows whose figures moved.
# - A handful of countries swapped adjacent table rows because the source
#   list was re-sorted; those rows' country-name cells (column A) are
#   corrected in place so each row keeps its own refreshed numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 25 de Junio de 2020 a las 21:53"

# Row 4: Estados Unidos - refreshed stats
$ws.Range("B4").Value = 2489121
$ws.Range("C4").Value = 26567
$ws.Range("D4").Value = 1044048
$ws.Range("E4").Value = 1318650
$ws.Range("G4").Value = 2142
$ws.Range("H4").Value = 126423

# Row 7: India - refreshed stats
$ws.Range("B7").Value = 491168
$ws.Range("C7").Value = 18183
$ws.Range("E7").Value = 190196

# Row 15: Alemania - refreshed stats
$ws.Range("B15").Value = 193465
$ws.Range("C15").Value = 211
$ws.Range("E15").Value = 7660
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 9005

# Row 19: Francia - refreshed stats
$ws.Range("D19").Value = 75351
$ws.Range("E19").Value = 56245
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = 29752

# Row 21: Sudafrica - refreshed stats
$ws.Range("B21").Value = 118375
$ws.Range("C21").Value = 6579
$ws.Range("D21").Value = 59974
$ws.Range("E21").Value = 56109
$ws.Range("G21").Value = 87
$ws.Range("H21").Value = 2292

# Row 30: Ecuador - refreshed stats
$ws.Range("B30").Value = 53156
$ws.Range("C30").Value = 1513
$ws.Range("D30").Value = 26097
$ws.Range("E30").Value = 22716
$ws.Range("G30").Value = 69
$ws.Range("H30").Value = 4343

# Row 75: corrected country name (Finlandia -> Uzbekistan); refreshed stats
$ws.Range("A75").Value = "Uzbekistan"
$ws.Range("B75").Value = 7177
$ws.Range("C75").Value = 276
$ws.Range("D75").Value = 4877
$ws.Range("E75").Value = 2280
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 20

# Row 76: corrected country name (Uzbekistan -> Finlandia); refreshed stats
$ws.Range("A76").Value = "Finlandia"
$ws.Range("B76").Value = 7172
$ws.Range("C76").Value = 5
$ws.Range("D76").Value = 6600
$ws.Range("E76").Value = 245
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 327

# Row 96: corrected country name (Tailandia -> Republica de Africa Central); refreshed stats
$ws.Range("A96").Value = "Republica de Africa Central"
$ws.Range("B96").Value = 3244
$ws.Range("C96").Value = 145
$ws.Range("D96").Value = 607
$ws.Range("E96").Value = 2597
$ws.Range("G96").Value = 2
$ws.Range("H96").Value = 40

# Row 97: corrected country name (Republica de Africa Central -> Tailandia); refreshed stats
$ws.Range("A97").Value = "Tailandia"
$ws.Range("B97").Value = 3158
$ws.Range("C97").Value = 1
$ws.Range("D97").Value = 3038
$ws.Range("E97").Value = 62
$ws.Range("H97").Value = 58

# Row 100: Costa Rica - refreshed stats
$ws.Range("B100").Value = 2684
$ws.Range("C100").Value = 169
$ws.Range("D100").Value = 1227
$ws.Range("E100").Value = 1445

# Row 129: corrected country name (Niger -> Yemen); refreshed stats
$ws.Range("A129").Value = "Yemen"
$ws.Range("B129").Value = 1076
$ws.Range("C129").Value = 61
$ws.Range("D129").Value = 386
$ws.Range("E129").Value = 402
$ws.Range("G129").Value = 14
$ws.Range("H129").Value = 288

# Row 130: corrected country name (Benin -> Niger); refreshed stats
$ws.Range("A130").Value = "Niger"
$ws.Range("B130").Value = 1056
$ws.Range("C130").Value = 5
$ws.Range("D130").Value = 917
$ws.Range("E130").Value = 72
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 67

# Row 131: corrected country name (Yemen -> Benin); refreshed stats
$ws.Range("A131").Value = "Benin"
$ws.Range("B131").Value = 1017
$ws.Range("C131").Value = 115
$ws.Range("D131").Value = 288
$ws.Range("E131").Value = 715
$ws.Range("G131").Value = 1
$ws.Range("H131").Value = 14

# Row 166: corrected country name (Guyana -> Angola); refreshed stats
$ws.Range("A166").Value = "Angola"
$ws.Range("B166").Value = 212
$ws.Range("C166").Value = 15
$ws.Range("D166").Value = 81
$ws.Range("E166").Value = 121
$ws.Range("H166").Value = 10

# Row 167: corrected country name (Angola -> Guyana); refreshed stats
$ws.Range("A167").Value = "Guyana"
$ws.Range("B167").Value = 209
$ws.Range("D167").Value = 107
$ws.Range("E167").Value = 90
$ws.Range("H167").Value = 12

# Row 202: corrected country name (Fiyi -> Dominica)
$ws.Range("A202").Value = "Dominica"

# Row 203: corrected country name (Dominica -> Fiyi)
$ws.Range("A203").Value = "Fiyi"

# Row 208: corrected country name (Groenlandia -> Islas Malvinas)
$ws.Range("A208").Value = "Islas Malvinas"

# Row 209: corrected country name (Islas Malvinas -> Groenlandia)
$ws.Range("A209").Value = "Groenlandia"

# Row 211: corrected country name (Seychelles -> Montserrat); refreshed stats
$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# Row 212: corrected country name (Montserrat -> Seychelles); refreshed stats
$ws.Range("A212").Value = "Seychelles"
$ws.Range("D212").Value = 11
$ws.Range("H212").Value = 0
